# Update the MIGO test output: new "NV CONTRATO" (P) and "NV PEDIDO" (Q)
# values reflecting the latest purchase-order execution for rows 2 and 3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = 4600244240
$ws.Range("Q2").Value = 4503342033
$ws.Range("P3").Value = 4600244241
$ws.Range("Q3").Value = 4503342035

# Leave the selection on the last updated cell, matching the author's
# on-screen state when the workbook was saved.
$ws.Range("Q3").Select()
